# The commit swaps the two player rows for Cagliari (Serie_A): the row that
# used to hold "Joseph Liteta" (row 14) now holds "Luca Mazzitelli", and the
# row that used to hold "Luca Mazzitelli" (row 15) now holds "Joseph Liteta".
# League/Team (columns A & B) and the trailing type/goalsPrevented columns
# (DK/DL) are identical between the two rows, so the net effect is simply
# swapping all the per-player stat columns (C through DL) between row 14 and
# row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row14Stats = $ws.Range("C14:DL14")
$row15Stats = $ws.Range("C15:DL15")

$buffer = $row14Stats.Value()
$row14Stats.Value = $row15Stats.Value()
$row15Stats.Value = $buffer
